# "modified BOM for new version box"
#
# The old BOM had two hardware rows that are no longer part of this
# revision of the box ("M1.4x4 Screw" + its order link, and "M3 Screw" +
# its order link), and a new part ("Slide switch") replaces them.
#
# Net effect on the sheet:
#   - the "M3 Screw" row is removed outright (everything below it shifts
#     up by one row)
#   - the row that used to hold "M1.4x4 Screw" now holds the new
#     "Slide switch" part (same row, same formatting - just new item name
#     + new order link)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "M3 Screw" line entirely - row 12 - shifting everything
# below (the Breadboard Jumper Wire Pack row, the spacer, the totals
# row and the Equipment/Usage table) up by one row.
$ws.Rows("12:12").Delete()

# Row 11 used to be "M1.4x4 Screw"; turn it into the new "Slide switch"
# part, keeping the row's existing number format / styling untouched.
$ws.Range("C11").Value = "http://www.amazon.com/SS12D00G3-Position-Panel-Vertical-Switch/dp/B007QAJUUS/ref=sr_1_7?s=industrial&ie=UTF8&qid=1406155714&sr=1-7&keywords=slide+switch"
$ws.Range("A11").Value = "Slide switch"

# Match the cursor position left behind in the saved workbook.
$ws.Range("A11").Select()
